$d = $word.ActiveDocument

$pairs = @(
    @("97×53=", "31×92="),
    @("25×99=", "50×89="),
    @("49×99=", "28×67="),
    @("46×32=", "59×63="),
    @("95×34=", "44×98="),
    @("51×88=", "49×52="),
    @("40×48=", "91×84="),
    @("75×78=", "79×61="),
    @("77×18=", "38×78="),
    @("92×79=", "91×87="),
    @("30×57=", "98×89="),
    @("60×72=", "66×59="),
    @("88×18=", "63×95="),
    @("43×36=", "47×72="),
    @("24×36=", "81×70="),
    @("65×19=", "16×83="),
    @("84×86=", "93×80="),
    @("30×14=", "46×99="),
    @("63×28=", "25×12="),
    @("96×34=", "21×72="),
    @("74×50=", "25×64="),
    @("56×63=", "93×19="),
    @("19×37=", "22×85="),
    @("56×82=", "45×44="),
    @("81×76=", "70×46=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
